$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.036.56'
$ws.Range("E2").Value = '  -4.72%  '
$ws.Range("D3").Value = '2.214.79'
$ws.Range("E3").Value = '  -6.17%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'315.01"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'100.46"
$ws.Range("E6").Value = '  -6.83%  '
$ws.Range("D7").Value = "'0.590"
$ws.Range("E7").Value = '  -6.15%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'0.565"
$ws.Range("E9").Value = '  -7.43%  '
$ws.Range("D10").Value = "'37.20"
$ws.Range("E10").Value = '  -9.01%  '
$ws.Range("D11").Value = "'54.11"
$ws.Range("E11").Value = '  -3.20%  '
$ws.Range("D12").Value = "'0.0832"
$ws.Range("E12").Value = '  -9.36%  '
$ws.Range("D13").Value = "'7.70"
$ws.Range("E13").Value = '  -9.50%  '
$ws.Range("D14").Value = "'0.108"
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = "'0.863"
$ws.Range("E15").Value = '  -11.59%  '
$ws.Range("D16").Value = '2.551.54'
$ws.Range("E16").Value = '  -6.14%  '
$ws.Range("D17").Value = "'14.29"
$ws.Range("E17").Value = '  -6.40%  '
$ws.Range("D18").Value = '2.214.04'
$ws.Range("E18").Value = '  -6.11%  '
$ws.Range("D19").Value = '42.935.86'
$ws.Range("E19").Value = '  -4.87%  '
$ws.Range("D20").Value = "'14.94"
$ws.Range("E20").Value = '  +5.83%  '
$ws.Range("D21").Value = '0.0₃0967'
$ws.Range("E21").Value = '  -8.94%  '
$ws.Range("D22").Value = "'6.50"
$ws.Range("E22").Value = '  -10.39%  '
$ws.Range("D23").Value = "'65.56"
$ws.Range("E23").Value = '  -10.42%  '
$ws.Range("D24").Value = "'3.14"
$ws.Range("E24").Value = '  -11.57%  '
$ws.Range("D25").Value = "'238.84"
$ws.Range("E25").Value = '  -7.55%  '
$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = '  -8.14%  '
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").Value = "'4.04"
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = "'10.03"
$ws.Range("E29").Value = '  -9.13%  '
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = '  -3.89%  '
$ws.Range("D31").Value = "'6.37"
$ws.Range("E31").Value = '  -11.90%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'20.65"
$ws.Range("E32").Value = '  -7.00%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0885"
$ws.Range("E33").Value = '  -9.64%  '
$ws.Range("D34").Value = "'34.74"
$ws.Range("E34").Value = '  -7.10%  '
$ws.Range("D35").Value = "'155.83"
$ws.Range("E35").Value = '  -6.64%  '
$ws.Range("D36").Value = "'2.79"
$ws.Range("E36").Value = '  -6.62%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = "'1.96"
$ws.Range("E37").Value = '  +9.51%  '
$ws.Range("D38").Value = "'0.122"
$ws.Range("E38").Value = '  -5.98%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = "'3.04"
$ws.Range("E39").Value = '  +5.37%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = "'4.42"
$ws.Range("E40").Value = '  -6.34%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = "'0.105"
$ws.Range("E41").Value = '  -10.30%  '
$ws.Range("E42").Value = '  -5.41%  '
$ws.Range("D43").Value = "'0.0327"
$ws.Range("E43").Value = '  -7.23%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '1.789.66'
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("D46").Value = "'12.20"
$ws.Range("E46").Value = '  -5.61%  '
$ws.Range("D47").Value = "'87.83"
$ws.Range("E47").Value = '  -10.41%  '
$ws.Range("D48").Value = "'0.206"
$ws.Range("E48").Value = '  -9.80%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = "'76.84"
$ws.Range("E49").Value = '  -9.80%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = "'5.35"
$ws.Range("E50").Value = '  -6.91%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = "'60.64"
$ws.Range("E51").Value = '  -13.00%  '
